# Junction_Flooding_373: refresh the 5 data rows with a newer measurement
# sample (custom accuracy / 1000-row dataset regeneration) and drop the
# now-unused 6th row, plus a handful of column width tweaks that came
# along with the regenerated sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (ColumnWidth setter adds ~0.8333 internally, so subtract 5/6 to land on target integer width)
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 6.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 6.166666666666667
$ws.Columns.Item(31).ColumnWidth = 4.166666666666667
$ws.Columns.Item(33).ColumnWidth = 7.166666666666667

# Update data rows 2-5 (new measurement data)
$ws.Range("A2").Value = 45054.50694444445
$ws.Range("B2").Value = 4.928
$ws.Range("C2").Value = 5.344
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 5.978
$ws.Range("F2").Value = 10.253
$ws.Range("G2").Value = 2.511
$ws.Range("H2").Value = 8.567
$ws.Range("I2").Value = 4.041
$ws.Range("J2").Value = 1.748
$ws.Range("K2").Value = 4.53
$ws.Range("L2").Value = 5.048
$ws.Range("M2").Value = 5.182
$ws.Range("N2").Value = 0.791
$ws.Range("O2").Value = 3.479
$ws.Range("P2").Value = 3.593
$ws.Range("Q2").Value = 1.503
$ws.Range("R2").Value = 1.042
$ws.Range("S2").Value = 0.417
$ws.Range("T2").Value = 43.098
$ws.Range("U2").Value = 7.946
$ws.Range("V2").Value = 4.36
$ws.Range("W2").Value = 6.645
$ws.Range("X2").Value = 2.748
$ws.Range("Y2").Value = 0.484
$ws.Range("Z2").Value = 2.765
$ws.Range("AA2").Value = 1.536
$ws.Range("AB2").Value = 3.068
$ws.Range("AC2").Value = 2.75
$ws.Range("AD2").Value = 5.612
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 5.035
$ws.Range("AG2").Value = 2.082
$ws.Range("AH2").Value = 3.31

$ws.Range("A3").Value = 45054.51388888889
$ws.Range("B3").Value = 24.827
$ws.Range("C3").Value = 19.296
$ws.Range("D3").Value = 0.635
$ws.Range("E3").Value = 51.803
$ws.Range("F3").Value = 45.261
$ws.Range("G3").Value = 18.771
$ws.Range("H3").Value = 68.378
$ws.Range("I3").Value = 29.065
$ws.Range("J3").Value = 13.175
$ws.Range("K3").Value = 20.517
$ws.Range("L3").Value = 21.87
$ws.Range("M3").Value = 23.076
$ws.Range("N3").Value = 6.049
$ws.Range("O3").Value = 19.216
$ws.Range("P3").Value = 26.846
$ws.Range("Q3").Value = 15.168
$ws.Range("R3").Value = 0.794
$ws.Range("S3").Value = 0.952
$ws.Range("T3").Value = 284.133
$ws.Range("U3").Value = 52.943
$ws.Range("V3").Value = 18.197
$ws.Range("W3").Value = 36.675
$ws.Range("X3").Value = 18.878
$ws.Range("Y3").Value = 2.527
$ws.Range("Z3").Value = 33.787
$ws.Range("AA3").Value = 15.013
$ws.Range("AB3").Value = 13.995
$ws.Range("AC3").Value = 16.276
$ws.Range("AD3").Value = 23.25
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 61.05
$ws.Range("AG3").Value = 10.319
$ws.Range("AH3").Value = 21.814

$ws.Range("A4").Value = 45054.52083333334
$ws.Range("B4").Value = 9.347
$ws.Range("C4").Value = 7.435
$ws.Range("D4").Value = 0.154
$ws.Range("E4").Value = 18.9
$ws.Range("F4").Value = 17.203
$ws.Range("G4").Value = 6.82
$ws.Range("H4").Value = 32.582
$ws.Range("I4").Value = 10.654
$ws.Range("J4").Value = 4.924
$ws.Range("K4").Value = 7.837
$ws.Range("L4").Value = 8.272
$ws.Range("M4").Value = 8.701000000000001
$ws.Range("N4").Value = 2.24
$ws.Range("O4").Value = 7.161
$ws.Range("P4").Value = 9.875999999999999
$ws.Range("Q4").Value = 5.431
$ws.Range("R4").Value = 0.485
$ws.Range("S4").Value = 0.443
$ws.Range("T4").Value = 101.183
$ws.Range("U4").Value = 19.77
$ws.Range("V4").Value = 6.894
$ws.Range("W4").Value = 13.968
$ws.Range("X4").Value = 7.042
$ws.Range("Y4").Value = 0.9429999999999999
$ws.Range("Z4").Value = 14.78
$ws.Range("AA4").Value = 5.405
$ws.Range("AB4").Value = 5.283
$ws.Range("AC4").Value = 6.144
$ws.Range("AD4").Value = 8.869
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 29.189
$ws.Range("AG4").Value = 3.931
$ws.Range("AH4").Value = 8.032

$ws.Range("A5").Value = 45054.52777777778
$ws.Range("B5").Value = 19.86
$ws.Range("C5").Value = 15.23
$ws.Range("D5").Value = 0.53
$ws.Range("E5").Value = 42.15
$ws.Range("F5").Value = 36.01
$ws.Range("G5").Value = 15.22
$ws.Range("H5").Value = 58.57
$ws.Range("I5").Value = 23.56
$ws.Range("J5").Value = 10.68
$ws.Range("K5").Value = 16.32
$ws.Range("L5").Value = 17.42
$ws.Range("M5").Value = 18.38
$ws.Range("N5").Value = 4.91
$ws.Range("O5").Value = 15.42
$ws.Range("P5").Value = 21.75
$ws.Range("Q5").Value = 12.4
$ws.Range("R5").Value = 0.48
$ws.Range("S5").Value = 0.72
$ws.Range("T5").Value = 227.12
$ws.Range("U5").Value = 42.75
$ws.Range("V5").Value = 14.44
$ws.Range("W5").Value = 29.4
$ws.Range("X5").Value = 15.28
$ws.Range("Y5").Value = 2.02
$ws.Range("Z5").Value = 28.25
$ws.Range("AA5").Value = 12.25
$ws.Range("AB5").Value = 11.15
$ws.Range("AC5").Value = 13.09
$ws.Range("AD5").Value = 18.44
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 52.45
$ws.Range("AG5").Value = 8.25
$ws.Range("AH5").Value = 17.63

# Remove old row 6 (dataset now has only 4 data rows instead of 5)
$ws.Rows.Item(6).Delete()
